# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-17 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-12-18 Wednesday", 2)

# Update the division-problem answers in the table. Cells are addressed
# directly by (row, column) so that duplicate answer strings elsewhere in
# the table are not ambiguously matched.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "68÷4=17, 0"
$t.Cell(1, 2).Range.Text = "85÷7=12, 1"
$t.Cell(1, 3).Range.Text = "50÷2=25, 0"
$t.Cell(1, 4).Range.Text = "73÷2=36, 1"
$t.Cell(1, 5).Range.Text = "92÷2=46, 0"

$t.Cell(5, 1).Range.Text = "67÷2=33, 1"
$t.Cell(5, 2).Range.Text = "33÷6=5, 3"
$t.Cell(5, 3).Range.Text = "37÷2=18, 1"
$t.Cell(5, 4).Range.Text = "57÷6=9, 3"
$t.Cell(5, 5).Range.Text = "34÷2=17, 0"

$t.Cell(9, 1).Range.Text = "93÷6=15, 3"
$t.Cell(9, 2).Range.Text = "35÷6=5, 5"
$t.Cell(9, 3).Range.Text = "63÷4=15, 3"
$t.Cell(9, 4).Range.Text = "77÷4=19, 1"
$t.Cell(9, 5).Range.Text = "81÷7=11, 4"

$t.Cell(13, 1).Range.Text = "80÷3=26, 2"
$t.Cell(13, 2).Range.Text = "66÷2=33, 0"
$t.Cell(13, 3).Range.Text = "92÷7=13, 1"
$t.Cell(13, 4).Range.Text = "30÷4=7, 2"
$t.Cell(13, 5).Range.Text = "90÷8=11, 2"

$t.Cell(17, 1).Range.Text = "64÷6=10, 4"
$t.Cell(17, 2).Range.Text = "43÷8=5, 3"
$t.Cell(17, 3).Range.Text = "67÷8=8, 3"
$t.Cell(17, 4).Range.Text = "66÷6=11, 0"
$t.Cell(17, 5).Range.Text = "57÷3=19, 0"

Write-Host "Done updating date and table cells."
